# Add a new "time_taken" metadata column (F) to the panel worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - copy formatting from the neighbouring header cell (E1)
# so the new header matches the existing bold/centered header style.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Data rows - plain text values, no special formatting.
$ws.Range("F2").Value = "2021-10-05 13:39:05.592193"
$ws.Range("F3").Value = "2021-10-05 13:39:05.592204"
